# Apply updated odds values to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$updates4 = @{
    "O4" = 1.3
    "P4" = 3.4
}
foreach ($addr in $updates4.Keys) {
    $ws.Range($addr).Value = $updates4[$addr]
}

# Row 6
$updates6 = @{
    "M6" = 1.07
    "N6" = 9
    "O6" = 1.36
    "P6" = 3
    "AG6" = 800
}
foreach ($addr in $updates6.Keys) {
    $ws.Range($addr).Value = $updates6[$addr]
}

# Row 8
$updates8 = @{
    "G8" = 2.12
    "H8" = 3.45
    "I8" = 3.15
    "J8" = 2.72
    "K8" = 2.12
    "L8" = 3.65
    "N8" = 8
    "P8" = 3.65
    "Q8" = 1.72
    "R8" = 2.05
    "S8" = 1.38
    "T8" = 2.8
    "U8" = 1.6
    "V8" = 2.2
    "W8" = 9
    "X8" = 11.5
    "Y8" = 8.5
    "Z8" = 21
    "AA8" = 15.5
    "AB8" = 22
    "AC8" = 8
    "AD8" = 6.8
    "AE8" = 12.5
    "AF8" = 50
    "AH8" = 11
    "AI8" = 17.5
    "AJ8" = 10.75
    "AK8" = 40
    "AL8" = 25
    "AM8" = 29
    "AN8" = 4.15
    "AO8" = 11
    "AP8" = 18.5
    "AQ8" = 40
    "AR8" = 70
    "AS8" = 200
    "AT8" = 2.8
    "AU8" = 6.9
    "AV8" = 60
    "AW8" = 5.1
    "AX8" = 17.5
    "AY8" = 23
    "AZ8" = 80
    "BA8" = 110
    "BB8" = 300
}
foreach ($addr in $updates8.Keys) {
    $ws.Range($addr).Value = $updates8[$addr]
}

# Row 9
$updates9 = @{
    "G9" = 2
    "I9" = 3.35
    "J9" = 2.6
    "L9" = 3.75
    "O9" = 1.24
    "T9" = 2.87
    "V9" = 2.15
    "W9" = 8.5
    "Y9" = 8.5
    "Z9" = 18
    "AA9" = 15
    "AB9" = 23
    "AH9" = 12
    "AI9" = 20
    "AJ9" = 11.5
    "AK9" = 45
    "AL9" = 27
    "AM9" = 30
    "AN9" = 4
    "AO9" = 10.25
    "AP9" = 18
    "AQ9" = 37
    "AR9" = 65
    "AT9" = 2.87
    "AU9" = 6.9
    "AW9" = 5.4
    "AX9" = 18
    "AY9" = 23
    "AZ9" = 90
    "BA9" = 110
    "BB9" = 250
}
foreach ($addr in $updates9.Keys) {
    $ws.Range($addr).Value = $updates9[$addr]
}

# Row 10
$updates10 = @{
    "L10" = 3.2
    "O10" = 1.29
    "Q10" = 1.88
    "W10" = 9
    "X10" = 13.5
    "AA10" = 21
    "AB10" = 28
    "AH10" = 8.75
    "AJ10" = 9.5
    "AL10" = 21
    "AM10" = 28
    "AO10" = 14
    "AP10" = 21
    "BB10" = 250
}
foreach ($addr in $updates10.Keys) {
    $ws.Range($addr).Value = $updates10[$addr]
}

# Row 11
$updates11 = @{
    "G11" = 1.42
    "H11" = 4.25
    "I11" = 7.3
    "J11" = 1.95
    "K11" = 2.25
    "L11" = 6.7
    "M11" = 1.05
    "N11" = 7.6
    "O11" = 1.27
    "P11" = 3.4
    "Q11" = 1.83
    "R11" = 1.91
    "S11" = 1.39
    "T11" = 2.77
    "X11" = 6.2
    "Z11" = 9
    "AC11" = 7.6
    "AH11" = 17
    "AI11" = 45
    "AJ11" = 23
    "AK11" = 175
    "AL11" = 90
    "AO11" = 6.6
    "AQ11" = 20
    "AT11" = 2.77
    "AV11" = 90
    "AW11" = 8.25
    "AX11" = 45
    "AY11" = 45
    "AZ11" = 350
}
foreach ($addr in $updates11.Keys) {
    $ws.Range($addr).Value = $updates11[$addr]
}
